$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (B1:F1) with new text/content, reordered ---
$ws.Range("F1").Value = "usser==null"
$ws.Range("E1").Value = 'user4=new User("Mahmud" , "s12     s123","memedoskim16@gmail.com") // space vo password'
$ws.Range("D1").Value = 'user3=new User("Mahmud" , "s123++s123","memedoskim16@gmail.com")//specail charcter in password'
$ws.Range("C1").Value = 'user2=new User("Mahmud" , "s123s123","mahmud.com")// no @ In email'
$ws.Range("B1").Value = 'user1=new User(null,"s123","aaa@test.com") // username = null'

# --- Apply wrap text formatting to header cells B1, D1, E1 ---
$ws.Range("B1").WrapText = $true
$ws.Range("D1").WrapText = $true
$ws.Range("E1").WrapText = $true

# --- Make sure C1 keeps its wrap-text + text-format style ---
$ws.Range("C1").WrapText = $true

# --- Column width changes ---
$ws.Columns.Item(4).ColumnWidth = 40.33203125
$ws.Columns.Item(5).ColumnWidth = 28.77734375

# --- Add a blank formatted cell D4 with wrap text ---
$ws.Range("D4").WrapText = $true

# --- Update selection to B1 ---
$ws.Range("B1").Select()
